# Insert two new daily price records for "Hortaliza, Feria Lagunitas de Puerto Montt - Lechuga"
# at rows 540:541, shifting the existing rows 540:579 down to 542:581.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above the current row 540 (pushes old 540..579 -> 542..581)
$ws.Rows("540:541").Insert()

# --- New row 540 ---
$ws.Range("A540").Value = 4
$ws.Range("B540").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C540").Value = "Los Lagos"
$ws.Range("D540").Value = 44746
$ws.Range("E540").Value = 10
$ws.Range("F540").Value = 100112033
$ws.Range("G540").Value = "Lechuga"
$ws.Range("H540").Value = "Escarola"
$ws.Range("I540").Value = "Primera"
$ws.Range("J540").Value = 200
$ws.Range("K540").Value = 15000
$ws.Range("L540").Value = 15000
$ws.Range("M540").Value = 15000
$ws.Range("N540").Value = "`$/caja 15 unidades"
$ws.Range("O540").Value = "Región de Coquimbo"
$ws.Range("P540").Value = 1000
$ws.Range("Q540").Value = 15
$ws.Range("R540").Value = "Hortaliza"

# --- New row 541 ---
$ws.Range("A541").Value = 4
$ws.Range("B541").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C541").Value = "Los Lagos"
$ws.Range("D541").Value = 44746
$ws.Range("E541").Value = 10
$ws.Range("F541").Value = 100112033
$ws.Range("G541").Value = "Lechuga"
$ws.Range("H541").Value = "Escarola"
$ws.Range("I541").Value = "Segunda"
$ws.Range("J541").Value = 100
$ws.Range("K541").Value = 13000
$ws.Range("L541").Value = 13000
$ws.Range("M541").Value = 13000
$ws.Range("N541").Value = "`$/caja 18 unidades"
$ws.Range("O541").Value = "Región de Coquimbo"
$ws.Range("P541").Value = 722
$ws.Range("Q541").Value = 18
$ws.Range("R541").Value = "Hortaliza"
